$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.256.24"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.789.84"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'226.14"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'0.550"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'32.29"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.294"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.0690"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'0.0946"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "2.049.38"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  -3.68%  "
$ws.Range("D14").Value = "1.792.23"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "34.236.89"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'68.04"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "0.0₃0803"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").Value = "'246.35"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "'10.98"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'4.19"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").Value = "'161.60"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'16.33"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("E33").Value = "  +3.76%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "1.441.01"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").Value = "'2.58"
$ws.Range("E36").Value = "  +9.91%  "
$ws.Range("D37").Value = "'0.667"
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0191"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.05"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "'81.84"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.39"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'14.08"
$ws.Range("E42").Value = "  +4.97%  "
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").Value = "'0.921"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'0.0520"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "'6.09"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "1.948.49"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "'105.53"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  -5.70%  "
